$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (column AC, the 29th column) - matches the text style of row 1
$headerCell = $ws.Range("AC1")
$headerCell.Value = "18-jul"
$headerCell.NumberFormat = "@"

# New values for column AC (rows 2-11), matching the centered-number style used in the rest of the table
$values = @(12, 17, 8, 9, 12, 15, 21, 12, 28, 24)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 29)
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
}

# Update the active selection to reflect the new last cell, as in the source workbook
[void]$ws.Range("AC12").Select()
